$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 6427
$ws.Range("L3").Value = 6928
$ws.Range("L4").Value = 1722
$ws.Range("L5").Value = 407
$ws.Range("L6").Value = 5682
$ws.Range("L7").Value = 21166

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L7").Value = 672
$ws.Range("L8").Value = 1401
$ws.Range("L15").Value = 177
$ws.Range("L19").Value = 579
$ws.Range("L20").Value = 538
$ws.Range("L24").Value = 65
$ws.Range("L27").Value = 183
$ws.Range("L29").Value = 1188
$ws.Range("L33").Value = 954
$ws.Range("L34").Value = 116
$ws.Range("L37").Value = 810
$ws.Range("L44").Value = 145
$ws.Range("L48").Value = 276
$ws.Range("L52").Value = 450
$ws.Range("L63").Value = 66
$ws.Range("L67").Value = 733
$ws.Range("L68").Value = 67
$ws.Range("L69").Value = 46
$ws.Range("L73").Value = 166
$ws.Range("L76").Value = 331
$ws.Range("L77").Value = 142
$ws.Range("L83").Value = 465
$ws.Range("L84").Value = 202
$ws.Range("L85").Value = 1049
$ws.Range("L90").Value = 226
$ws.Range("L91").Value = 283
$ws.Range("L92").Value = 66
$ws.Range("L99").Value = 367
$ws.Range("L101").Value = 21166

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L2").Value = 235
$ws.Range("L7").Value = 672

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L3").Value = 437
$ws.Range("L6").Value = 215
$ws.Range("L7").Value = 1049

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L2").Value = 144
$ws.Range("L7").Value = 450

$ws = $wb.Worksheets.Item("Norwood Park")
$ws.Range("L2").Value = 20
$ws.Range("L7").Value = 46

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L4").Value = 97
$ws.Range("L7").Value = 1401

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L3").Value = 186
$ws.Range("L6").Value = 105
$ws.Range("L7").Value = 465

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L2").Value = 262
$ws.Range("L3").Value = 336
$ws.Range("L7").Value = 954

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L3").Value = 286
$ws.Range("L4").Value = 47
$ws.Range("L6").Value = 209
$ws.Range("L7").Value = 810

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L2").Value = 113
$ws.Range("L7").Value = 367

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L2").Value = 209
$ws.Range("L3").Value = 285
$ws.Range("L4").Value = 49
$ws.Range("L7").Value = 733

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("L5").Value = 4
$ws.Range("L7").Value = 202

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 361
$ws.Range("L3").Value = 457
$ws.Range("L7").Value = 1188

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L6").Value = 112
$ws.Range("L7").Value = 276

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L3").Value = 177
$ws.Range("L6").Value = 158
$ws.Range("L7").Value = 579

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("L3").Value = 41
$ws.Range("L7").Value = 145

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L6").Value = 148
$ws.Range("L7").Value = 331

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("L3").Value = 22
$ws.Range("L7").Value = 65

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L2").Value = 99
$ws.Range("L7").Value = 283

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L3").Value = 188
$ws.Range("L4").Value = 53
$ws.Range("L7").Value = 538

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("L2").Value = 41
$ws.Range("L7").Value = 116

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("L3").Value = 57
$ws.Range("L7").Value = 177

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("L6").Value = 37
$ws.Range("L7").Value = 166

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("L3").Value = 13
$ws.Range("L7").Value = 66

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("L6").Value = 58
$ws.Range("L7").Value = 183

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("L2").Value = 76
$ws.Range("L3").Value = 64
$ws.Range("L7").Value = 226

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("L2").Value = 23
$ws.Range("L7").Value = 67

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("L3").Value = 45
$ws.Range("L7").Value = 142
